# feature multithread for saving data
#
# The user-import worksheet gains four more imported rows (5-8), row 3 is
# completed with its missing rownumber/username, and row 4's email is
# corrected (it was importing with row 3's address). Hyperlinks are added
# for every new email cell, reusing the sheet's existing "Hyperlink" cell
# style. The writes below are ordered to match the concurrent/multithreaded
# worker completion order used by the importer (row 5 finishes, then row 6
# starts but stalls, rows 7 and 8 race ahead and finish, row 3 and then
# row 6 finally complete, row 4 is patched, and the last worker's roletype
# cell lands last).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 5 (worker finishes quickly, end-to-end) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "test4"
$ws.Range("D5").Value = "test4@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:test4@gmail.com")
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("C5").Value = "pass12313123"
$ws.Range("E5").Value = "MANAGER, USER"
$ws.Range("F5").Value = "MANAGER, USER"

# --- row 6 (worker writes rownumber/username, then stalls) ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "test5"

# --- row 7 (races ahead of row 6's stalled worker) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "test6"
$ws.Range("D7").Value = "test6@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:test6@gmail.com")
$ws.Range("D7").Style = "Hyperlink"
$ws.Range("C7").Value = "pass12313123"
$ws.Range("E7").Value = "MANAGER, USER"
$ws.Range("F7").Value = "MANAGER, USER"

# --- row 8 (also races ahead; roletype cell written later, see below) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "test7"
$ws.Range("D8").Value = "test7@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:test7@gmail.com")
$ws.Range("D8").Style = "Hyperlink"
$ws.Range("C8").Value = "pass12313123"
$ws.Range("E8").Value = "MANAGER, USER"

# --- row 3 gets its missing rownumber/username filled in ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "test2"

# --- row 6's worker resumes and finishes ---
$ws.Range("D6").Value = "test5@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:test5@gmail.com")
$ws.Range("D6").Style = "Hyperlink"
$ws.Range("C6").Value = "pass12313123"
$ws.Range("E6").Value = "MANAGER, USER"
$ws.Range("F6").Value = "MANAGER, USER"

# --- row 4's email gets corrected to its own address, rolecode backfilled ---
$ws.Range("D4").Value = "test3@gmail.com"
$ws.Range("E4").Value = "ADMIN, MANAGER"

# --- row 8's roletype is the last cell to land ---
$ws.Range("F8").Value = "USER, USER"

# Leave the selection where the importer's last active cell ended up.
$null = $ws.Range("C13").Select()
